# Remove testers who didn't end up testing:
#   row 20 - St. Norbert / Wisconsin
#   row 18 - Wellesley / Massachusetts
#   row 17 - Longwood / Virginia (duplicate entry)
# Delete from bottom to top so earlier row numbers stay valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(20).Delete()
$ws.Rows(18).Delete()
$ws.Rows(17).Delete()

$ws.Range("D9").Select()
